# Apply the "changed model paths file to better default" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filenames")

# Update the path cells to use relative paths instead of absolute macOS paths.
$ws.Range("B2").Value = "model_parameters/"
$ws.Range("B3").Value = "input_data/"
$ws.Range("B4").Value = "output_data/"
$ws.Range("B5").Value = "output_figures/"

# Update the default parameter table filename to the example table.
$ws.Range("B6").Value = "Road_dust_parameter_table_example.xlsx"

# Update the default output data filename and add an alternative in C8.
$ws.Range("B8").Value = "example_output.xlsx"
$ws.Range("C8").Value = "output.txt"
